$d = $word.ActiveDocument

# --- Paragraph 1: "Verificare funzionamento ... Stampa Ord Frut Loc" ---
# Replace the whole paragraph content with a green-highlighted version where
# "Frut" and "Loc" each get their own run + spell-check proofErr markers,
# and drop the trailing _GoBack bookmark (it moves to the next paragraph).
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Verificare funzionamento*") {
        $p1 = $p
        break
    }
}

$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00E90BFB" w:rsidRDefault="00E90BFB" w:rsidP="002937FF">
<w:pPr>
<w:pStyle w:val="Paragrafoelenco"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr>
<w:spacing w:line="256" w:lineRule="auto"/>
<w:rPr><w:highlight w:val="green"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Verificare funzionamento &#8220;Stampa </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Ord</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Frut</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Loc</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>&#8221;</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$p1.Range.InsertXML($xml1)

# --- Paragraph 2: "Aggiungere alla ... Stampa Ord Frut Cant Excel ..." ---
# Add green highlight to the paragraph mark + the run, and move the
# _GoBack bookmark here (at the very start of the paragraph).
$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Aggiungere alla*") {
        $p2 = $p
        break
    }
}

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00E90BFB" w:rsidRDefault="00E90BFB" w:rsidP="000E0AED">
<w:pPr>
<w:pStyle w:val="Paragrafoelenco"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr>
<w:spacing w:line="256" w:lineRule="auto"/>
<w:rPr><w:highlight w:val="green"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Aggiungere alla &#8220;Stampa Ord Frut Cant Excel&#8221; la lista dei frutti (Non appartenenti ad un gruppo)</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$p2.Range.InsertXML($xml2)

# Re-find paragraph 2 after the replace and drop the bookmark at its start.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Aggiungere alla*") {
        $startRng = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $startRng)
        break
    }
}
